# Update the deliverables checklist:
#  - "Presentation slide deck" moves into row 3 (was row 4), marked fully done.
#  - "Assessment report" moves into row 4 (was row 3), status note refreshed.
# This mirrors swapping the two rows' contents (values + row height) in place,
# along with a couple of value/text tweaks called out in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) row 3 and row 4 contents so we can swap them.
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$e3 = $ws.Range("E3").Value2
$height3 = $ws.Rows.Item(3).RowHeight

$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$c4 = $ws.Range("C4").Value2
$e4 = $ws.Range("E4").Value2
$height4 = $ws.Rows.Item(4).RowHeight

# New row 4: the "Assessment report" entry, with an updated status note.
$ws.Range("A4").Value = $a3
$ws.Range("B4").Value = $b3
$ws.Range("C4").Value = $c3
$ws.Range("D4").Value = "the file exists and is almost done"
$ws.Range("E4").Value = $e3
$ws.Rows.Item(4).RowHeight = $height3

# New row 3: the "Presentation slide deck" entry, now complete.
$ws.Range("A3").Value = $a4
$ws.Range("B3").Value = $b4
$ws.Range("C3").Value = $c4
$ws.Range("D3").Value = "slides are done"
$ws.Range("E3").Value = 2
$ws.Rows.Item(3).RowHeight = $height4

# Reflect where the editor's cursor ended up: active cell D4.
$ws.Range("D4").Select()
